# Actualización automática 2025-07-03 11:34:50
# Update PRESUPUESTO (column G) values on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    2  = 2000
    3  = 3500
    4  = 2000
    6  = 2000
    9  = 3000
    10 = 10000
    12 = 3500
    14 = 1000
    16 = 2000
    19 = 500
    22 = 1000
    26 = 3000
    27 = 2000
    38 = 0
    39 = 1000
    44 = 1000
    45 = 4000
    49 = 0
    51 = 2000
    53 = 1000
    55 = 0
    58 = 54500
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
